$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new column before N (shifts N..X to O..Y) and give it the
# header "利率種類". The new cell inherits column N's old header style (s=6),
# which is exactly what's needed for N1 going forward.
$ws.Columns("N").Insert()
$ws.Range("N1").Value = "利率種類"

# Match the target column widths for the new column (N) and the one
# immediately left of it (M), which was narrowed in this edit.
# (stored xlsx <col> width = ColumnWidth + 5/7)
$ws.Columns("M").ColumnWidth = 10.44140625 - 5/7
$ws.Columns("N").ColumnWidth = 9.5546875 - 5/7

# Expand the hidden AutoFilter defined name by one column (P -> Q) to
# account for the newly inserted column.
foreach ($n in $wb.Names) {
    if ($n.Name -like "*_FilterDatabase*") {
        $n.RefersTo = "=正常件!`$A`$1:`$Q`$1"
    }
}

# Restore the active selection to the cell the author ended up on.
$ws.Range("Q8").Select()
